$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values (A2:D4)
$ws.Range("A2").Value = 47.425127192557575
$ws.Range("B2").Value = 0.99334575871763475
$ws.Range("C2").Value = 0.099167196384885331
$ws.Range("D2").Value = -0.058566806288615825

$ws.Range("A3").Value = -9.740285675559285
$ws.Range("B3").Value = 0.11267508584553425
$ws.Range("C3").Value = -0.73150735677226075
$ws.Range("D3").Value = 0.67245915267602829

$ws.Range("A4").Value = -80.775434919996215
$ws.Range("B4").Value = 0.023843839191458899
$ws.Range("C4").Value = -0.67458346714785578
$ws.Range("D4").Value = -0.73781340268620055

# Update the selection to B2:D4 with active cell B2
$ws.Range("B2:D4").Select()
